$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 425, shifting existing rows 425-485 down to 426-486.
$ws.Rows(425).Insert()

# Populate the newly inserted row 425 with the new weekly record.
$ws.Cells.Item(425, 1).Value = 4
$ws.Cells.Item(425, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(425, 3).Value = "Los Lagos"
$ws.Cells.Item(425, 4).Value = 45154
$ws.Cells.Item(425, 5).Value = 10
$ws.Cells.Item(425, 6).Value = 100112043
$ws.Cells.Item(425, 7).Value = "Pepino ensalada"
$ws.Cells.Item(425, 8).Value = "Sin especificar"
$ws.Cells.Item(425, 9).Value = "Primera"
$ws.Cells.Item(425, 10).Value = 120
$ws.Cells.Item(425, 11).Value = 15000
$ws.Cells.Item(425, 12).Value = 15000
$ws.Cells.Item(425, 13).Value = 15000
$ws.Cells.Item(425, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(425, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(425, 16).Value = 250
$ws.Cells.Item(425, 17).Value = 60
$ws.Cells.Item(425, 18).Value = "Hortaliza"
